# Updated: po 09. 08. 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows (AgTests / AgPosit columns) ---
$ws.Range("F369").Value = 235636

$ws.Range("F377").Value = 177042

$ws.Range("F391").Value = 178271

$ws.Range("F407").Value = 158562

$ws.Range("F424").Value = 266818

$ws.Range("F453").Value = 70391

$ws.Range("F474").Value = 46052

$ws.Range("F475").Value = 36673

$ws.Range("F477").Value = 37135

$ws.Range("F483").Value = 65645

$ws.Range("F499").Value = 11305

$ws.Range("F510").Value = 7777

$ws.Range("F511").Value = 6742

$ws.Range("F512").Value = 8389

$ws.Range("F513").Value = 10239

$ws.Range("F514").Value = 6857

$ws.Range("F515").Value = 4943
$ws.Range("G515").Value = 15

$ws.Range("F516").Value = 9187
$ws.Range("G516").Value = 32

$ws.Range("F517").Value = 6630
$ws.Range("G517").Value = 14

$ws.Range("F518").Value = 6936
$ws.Range("G518").Value = 9

# --- New rows appended at the bottom (519-522) ---
$newRows = @(
    @{ Row = 519; A = 44413; B = 392963; C = 7265;  D = 65; E = 12541; F = 7681; G = 18 },
    @{ Row = 520; A = 44414; B = 393007; C = 10312; D = 44; E = 12541; F = 8993; G = 22 },
    @{ Row = 521; A = 44415; B = 393040; C = 4290;  D = 33; E = 12541; F = 5594; G = 20 },
    @{ Row = 522; A = 44416; B = 393059; C = 1881;  D = 19; E = 12541; F = 3329; G = 10 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "yyyy-mm-dd"
    $cellA.Value = $r.A

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
}
